$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column K data (header year 2020, value 173)
$ws.Range("K4").Value = 2020
$ws.Range("K5").Value = 173

# Copy style from column J to column K for rows 3,4,5
$ws.Range("J3").Copy() | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("J4").Copy() | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null

$ws.Range("J5").Copy() | Out-Null
$ws.Range("K5").PasteSpecial(-4122) | Out-Null

$ws.Range("K4").Value = 2020
$ws.Range("K5").Value = 173

# Update selection on the sheet view
$ws.Range("I18").Select() | Out-Null
